$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1 -- copy the formatting (bold/border/alignment)
# from the neighboring header cell G1 so it reuses the same style record.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill H2:H12 with 0 for each data row (plain, unstyled numeric cells).
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
